$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Append-After($anchorText, $newText) {
    # Find `anchorText`, collapse the found range to its end, then insert
    # `newText` right after it (used to splice in brand-new sentences).
    $rng = $d.Content
    $ok = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $rng.Collapse(0)
        $rng.InsertAfter($newText)
    }
}

# --- Title ---
Replace-Text "Unraveling the Mysteries of Space: A Journey Through the Galaxy" "Exploring the Convergence of Art and Science"

# --- Author line: collapse "Dr" + "." + " Emily Carter" into a single new name ---
Replace-Text "Dr. Emily Carter" "Carissa Fernandez"

# --- Email line ---
Replace-Text "EmilyCarterPhD@cosmosresearch" "fernandezcarissa07@gmail"
Replace-Text "edu" "com"

# --- Body paragraph sentence-by-sentence replacement ---
Replace-Text "The vast expanse of the cosmos has captivated humanity for centuries, inspiring awe and wonder" "Art and science, often perceived as disparate disciplines, share an intrinsic connection that weaves together creativity and rationality"

Replace-Text " From the earliest astronomers gazing up at the night sky to the modern era of space exploration, we have embarked on an ongoing journey to understand the intricacies of the universe" " They both stem from a profound curiosity and an insatiable desire to understand and express the intricacies of the world around us"

Replace-Text " This exploration has led to profound insights into the nature of our place in the cosmos and the fundamental laws that govern the universe" " In the realm of art, we find emotions and imagination taking center stage, while in the domain of science, logic and reason lead the way"

Replace-Text "We have witnessed the birth and death of stars, the formation of galaxies, and the enigmatic phenomena of black holes" "The convergence of art and science is evident in the ways that artistic expression can illuminate scientific concepts, making them more accessible and engaging"

Replace-Text " We have discovered planets orbiting distant suns, raising questions about the potential for life beyond Earth" " Through paintings, sculptures, music, and literature, artists have the unique ability to translate complex scientific phenomena into forms that resonate with our senses and emotions, fostering a deeper understanding and appreciation of the natural world"

Replace-Text " The mysteries of space continue to beckon us, fueling our insatiable curiosity and driving our quest for knowledge" " Conversely, science provides art with a rich tapestry of inspiration, offering artists a boundless realm of forms, colors, and patterns to draw upon"

Replace-Text "With each new discovery, we deepen our understanding of the universe and our place within it" "Furthermore, both art and science share a common goal: to communicate ideas and inspire thought"

Replace-Text " We unravel the secrets of cosmic evolution, unraveling the history of the universe and tracing its trajectory into the future" " Artists strive to convey their perspectives, emotions, and experiences through their works, while scientists aim to share their findings and insights with the world"

Replace-Text " Our journey through the galaxy is a testament to humanity's enduring quest for knowledge and our unwavering fascination with the boundless mysteries of space" " Both disciplines rely on effective communication to engage their audiences, whether it be through the evocative power of imagery or the persuasive force of logical argument"

# --- New sentences appended within the body paragraph ---
Append-After "logic and reason lead the way." " However, upon closer examination, the boundaries between these two seemingly contrasting realms begin to blur, revealing a captivating interplay that has shaped human thought and culture throughout history"

Append-After "a boundless realm of forms, colors, and patterns to draw upon." " The breathtaking beauty of a starry night sky, the intricate structure of a flower, or the rhythmic pulse of a heartbeat can all serve as muses, igniting the creative spark in an artist's mind"

Append-After "the persuasive force of logical argument." " The ability to effectively communicate complex concepts is essential for both artists and scientists, and it is through this shared purpose that they find common ground"

# --- Summary heading paragraph stays "Summary" (unchanged) ---

# --- Summary body paragraph ---
Replace-Text "Our exploration of space has yielded remarkable insights into the nature and history of the universe" "The convergence of art and science is a testament to the multifaceted nature of human understanding"

Replace-Text " We have witnessed celestial wonders, from star formations to black holes, and discovered planets beyond our solar system" " Through their unique perspectives, artists and scientists complement each other, offering a holistic approach to comprehending the universe"

Replace-Text " The pursuit of space exploration continues to drive our quest for knowledge, captivating humanity with its profound implications for our understanding of the universe and our place within it" " Art illuminates the emotional and intuitive dimensions of existence, while science provides a framework for rational inquiry and empirical evidence"

Append-After "science provides a framework for rational inquiry and empirical evidence." " Together, they create a dynamic interplay that enriches our understanding of the world and fuels the progress of human thought and culture"

# --- New trailing empty paragraph before the section break ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

Write-Output "edit complete"
